# Fix typo in Spanish axis/header labels ("serorreversion" -> "seroreversion")
# on the summary row (row 8) of the sheet.
#
# Commit message: "cambio de de debut sexual y nombres de ejes en español"
#   (fix of sexual debut and axis names in Spanish)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the double-"rr" typo in the two "Tasa de serorreversion..." headers
# located in row 8 (columns F and G).
$ws.Range("F8").Value = "Tasa de seroreversion"
$ws.Range("G8").Value = "Tasa de seroreversión Rhat"
